$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = 1

# Update row 18 values
$ws.Range("C18").Value = 19
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 16

# Update selection to D18
$ws.Range("D18").Select()
